$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E) for rows 16-71 so periods run in ascending
# chronological order (1612 .. 2107) instead of descending (2107 .. 1612).
$periodos = @(
  "1612","1701","1702","1703","1704","1705","1706","1707","1708","1709",
  "1710","1711","1712","1801","1802","1803","1804","1805","1806","1807",
  "1808","1809","1810","1811","1812","1901","1902","1903","1904","1905",
  "1906","1907","1908","1909","1910","1911","1912","2001","2002","2003",
  "2004","2005","2006","2007","2008","2009","2010","2011","2012","2101",
  "2102","2103","2104","2105","2106","2107"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
  $row = $startRow + $i
  $ws.Range("E$row").Value = $periodos[$i]
}

# The "Valor Mora" (F) amounts follow the worker's period, so they move
# together with the period that used to sit in that row: period 1612 (now
# row 16) carries 120000 and period 2107 (now row 71) carries 100000.
$ws.Range("F16").Value = 120000
$ws.Range("F71").Value = 100000
